$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 452) holds a date serial value that was
# incremented by one day (45180 -> 45181, i.e. 2023-09-11 -> 2023-09-12).
$ws.Range("C2:C452").Value = 45181
